$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2823.8
$ws.Range("I98").Value = 451.0909
$ws.Range("J98").Value = 9348.75
$ws.Range("K98").Value = 451.0909
$ws.Range("L98").Value = 9348.75
$ws.Range("M98").Value = 1046.9091
$ws.Range("N98").Value = -12344.75

$ws.Range("H116").Value = 2740
$ws.Range("I116").Value = 2425
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 2425
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 1017
$ws.Range("N116").Value = -10884

$ws.Range("H122").Value = 2823.8
$ws.Range("I122").Value = 451.0909
$ws.Range("J122").Value = 9348.75
$ws.Range("K122").Value = 1353.2727
$ws.Range("L122").Value = 28046.25
$ws.Range("M122").Value = 1096.7273
$ws.Range("N122").Value = -32946.25

$ws.Range("H141").Value = 3157.5454
$ws.Range("I141").Value = 1882
$ws.Range("K141").Value = 5646
$ws.Range("M141").Value = -466

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 708
$ws.Range("I2").Value = 708
$ws.Range("K2").Value = 708
$ws.Range("M2").Value = -595

$ws.Range("H32").Value = 9144.550999999999
$ws.Range("I32").Value = 9482.23
$ws.Range("K32").Value = 9482.23
$ws.Range("M32").Value = -9195.23

$ws.Range("H74").Value = 1312.6757
$ws.Range("I74").Value = 1341.96
$ws.Range("J74").Value = 1251.6666
$ws.Range("K74").Value = 1341.96
$ws.Range("L74").Value = 1251.6666
$ws.Range("M74").Value = -467.96
$ws.Range("N74").Value = -2999.6666

$ws.Range("H77").Value = 1312.6757
$ws.Range("I77").Value = 1341.96
$ws.Range("J77").Value = 1251.6666
$ws.Range("K77").Value = 6709.8
$ws.Range("L77").Value = 6258.333000000001
$ws.Range("M77").Value = -2341.8
$ws.Range("N77").Value = -14994.333

$ws.Range("H102").Value = 2066.6667
$ws.Range("I102").Value = 1933.3334
$ws.Range("J102").Value = 2333.3333
$ws.Range("K102").Value = 1933.3334
$ws.Range("L102").Value = 2333.3333
$ws.Range("M102").Value = -311.3334
$ws.Range("N102").Value = -5577.3333

$ws.Range("H110").Value = 2350
$ws.Range("I110").Value = 2166.6667
$ws.Range("K110").Value = 2166.6667
$ws.Range("M110").Value = -121.6667000000002

$ws.Range("H116").Value = 708
$ws.Range("I116").Value = 708
$ws.Range("K116").Value = 708
$ws.Range("M116").Value = 1586

$ws.Range("H132").Value = 8065971
$ws.Range("I132").Value = 10417628
$ws.Range("J132").Value = 3147.1428
$ws.Range("K132").Value = 31252884
$ws.Range("L132").Value = 9441.428400000001
$ws.Range("M132").Value = -31250354
$ws.Range("N132").Value = -14501.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 708
$ws.Range("I3").Value = 708
$ws.Range("K3").Value = 708
$ws.Range("M3").Value = -594

$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -502
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 3560.6
$ws.Range("I105").Value = 1687.2142
$ws.Range("K105").Value = 1687.2142
$ws.Range("M105").Value = 59.78580000000011

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1583.8462
$ws.Range("I105").Value = 1570
$ws.Range("K105").Value = 1570
$ws.Range("M105").Value = 177

$ws.Range("H134").Value = 823692.2
$ws.Range("I134").Value = 2865.9092
$ws.Range("J134").Value = 3403432
$ws.Range("K134").Value = 8597.7276
$ws.Range("L134").Value = 10210296
$ws.Range("M134").Value = -6062.7276
$ws.Range("N134").Value = -10215366

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 775.38464
$ws.Range("J98").Value = 761.8182
$ws.Range("L98").Value = 2285.4546
$ws.Range("N98").Value = -5281.4546

$ws.Range("H114").Value = 1279.2632
$ws.Range("I114").Value = 301.57144
$ws.Range("J114").Value = 1849.5834
$ws.Range("K114").Value = 904.71432
$ws.Range("L114").Value = 5548.7502
$ws.Range("M114").Value = 2349.28568
$ws.Range("N114").Value = -12056.7502

$ws.Range("H131").Value = 828.59
$ws.Range("I131").Value = 419.9091
$ws.Range("J131").Value = 879.10114
$ws.Range("K131").Value = 1259.7273
$ws.Range("L131").Value = 2637.30342
$ws.Range("M131").Value = 3780.2727
$ws.Range("N131").Value = -12717.30342

$ws.Range("H132").Value = 1156.4517
$ws.Range("I132").Value = 725.55554
$ws.Range("J132").Value = 1753.0769
$ws.Range("K132").Value = 6529.99986
$ws.Range("L132").Value = 15777.6921
$ws.Range("M132").Value = -3999.99986
$ws.Range("N132").Value = -20837.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11114585
$ws.Range("I80").Value = 25643560
$ws.Range("J80").Value = 4192.294
$ws.Range("K80").Value = 25643560
$ws.Range("L80").Value = 4192.294
$ws.Range("M80").Value = -25642562
$ws.Range("N80").Value = -6188.294

$ws.Range("H83").Value = 11114585
$ws.Range("I83").Value = 25643560
$ws.Range("J83").Value = 4192.294
$ws.Range("K83").Value = 128217800
$ws.Range("L83").Value = 20961.47
$ws.Range("M83").Value = -128212808
$ws.Range("N83").Value = -30945.47

$ws.Range("H113").Value = 91968.63
$ws.Range("I113").Value = 125903.875
$ws.Range("J113").Value = 1474.6666
$ws.Range("K113").Value = 125903.875
$ws.Range("L113").Value = 1474.6666
$ws.Range("M113").Value = -123733.875
$ws.Range("N113").Value = -5814.6666

$ws.Range("H126").Value = 4746
$ws.Range("I126").Value = 2082.4
$ws.Range("J126").Value = 5633.8667
$ws.Range("K126").Value = 6247.200000000001
$ws.Range("L126").Value = 16901.6001
$ws.Range("M126").Value = -3777.200000000001
$ws.Range("N126").Value = -21841.6001

$ws.Range("H132").Value = 2936.7354
$ws.Range("I132").Value = 2747.9583
$ws.Range("J132").Value = 3389.8
$ws.Range("K132").Value = 8243.874899999999
$ws.Range("L132").Value = 10169.4
$ws.Range("M132").Value = -5713.874899999999
$ws.Range("N132").Value = -15229.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1107
$ws.Range("I16").Value = 1121.7858
$ws.Range("K16").Value = 1121.7858
$ws.Range("M16").Value = -951.7858000000001

$ws.Range("H40").Value = 8931.125
$ws.Range("I40").Value = 11999.6
$ws.Range("K40").Value = 11999.6
$ws.Range("M40").Value = -11863.6

$ws.Range("H46").Value = 965.087
$ws.Range("I46").Value = 533.3333
$ws.Range("J46").Value = 1029.85
$ws.Range("K46").Value = 533.3333
$ws.Range("L46").Value = 1029.85
$ws.Range("M46").Value = -345.3333
$ws.Range("N46").Value = -1405.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2481.4
$ws.Range("I107").Value = 2976
$ws.Range("J107").Value = 503
$ws.Range("K107").Value = 8928
$ws.Range("L107").Value = 1509
$ws.Range("M107").Value = -7008
$ws.Range("N107").Value = -5349

$ws.Range("H126").Value = 2467.037
$ws.Range("I126").Value = 1742.4
$ws.Range("J126").Value = 4537.4287
$ws.Range("K126").Value = 5227.200000000001
$ws.Range("L126").Value = 13612.2861
$ws.Range("M126").Value = -2757.200000000001
$ws.Range("N126").Value = -18552.2861

$ws.Range("H136").Value = 835.5185
$ws.Range("I136").Value = 720.3077
$ws.Range("J136").Value = 1135.0667
$ws.Range("K136").Value = 2160.9231
$ws.Range("L136").Value = 3405.2001
$ws.Range("M136").Value = 389.0769
$ws.Range("N136").Value = -8505.2001
